# Update the "dSF" (F) column values for a number of rows.
# These values represent re-pulled / recalculated data per the commit
# message ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -3
    6  = -2
    7  = 0
    10 = -7
    12 = -2
    16 = 4
    18 = 2
    26 = 4
    32 = -1
    34 = 2
    37 = 1
    38 = -1
    42 = 0
    48 = 0
    51 = 0
    54 = 2
    58 = 4
    60 = -1
    61 = -2
    66 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
